$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the old "Sample Year" column (E).
$ws.Range("E1:E2").EntireColumn.Delete()

# 2) Insert a new column for "Schulman Shift" right after "Sample Year Growth" (now E).
$ws.Range("F1:F2").EntireColumn.Insert()

# 3) Insert two new columns before the old "Sample note" column (now N) for
#    "Sample height (m)" and "Sample azimuth ()".
$ws.Range("N1:N2").EntireColumn.Insert()
$ws.Range("N1:N2").EntireColumn.Insert()

# --- Header row ---
$ws.Range("F1").Value = "Schulman Shift"
$ws.Range("N1").Value = "Sample height (m)"
$ws.Range("O1").Value = "Sample azimuth ()"

# --- Data row ---
$ws.Range("F2").Value = $false
$ws.Range("F2").NumberFormat = '"TRUE";"TRUE";"FALSE"'
$ws.Range("N2").Value = 1.5
$ws.Range("O2").Value = "NA"

$ws.Range("F6").Select()
